$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2062937062937063
$ws.Range("C2").Value = 0.5734265734265734
$ws.Range("J2").Value = 0.02447552447552448
$ws.Range("P2").Value = 0.1363636363636364
$ws.Range("S2").Value = 0.05944055944055944
$ws.Range("C3").Value = 0.06214689265536723
$ws.Range("J3").Value = 0.03954802259887006
$ws.Range("P3").Value = 0.7175141242937854
$ws.Range("S3").Value = 0.1807909604519774
$ws.Range("J4").Value = 0.06382978723404255
$ws.Range("P4").Value = 0.5319148936170213
$ws.Range("S4").Value = 0.4042553191489361
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.05294117647058823
$ws.Range("D6").Value = 0.03529411764705882
$ws.Range("F6").Value = 0.02941176470588235
$ws.Range("J6").Value = 0.2058823529411765
$ws.Range("O6").Value = 0.03529411764705882
$ws.Range("Q6").Value = 0.1352941176470588
$ws.Range("R6").Value = 0.07058823529411765
$ws.Range("S6").Value = 0.4352941176470588
$ws.Range("B7").Value = 0.1073825503355705
$ws.Range("D7").Value = 0.01342281879194631
$ws.Range("E7").Value = 0.01342281879194631
$ws.Range("F7").Value = 0.04026845637583892
$ws.Range("J7").Value = 0.1208053691275168
$ws.Range("O7").Value = 0.02684563758389262
$ws.Range("Q7").Value = 0.2080536912751678
$ws.Range("R7").Value = 0.09395973154362416
$ws.Range("S7").Value = 0.3758389261744967
$ws.Range("B8").Value = 0.1002044989775051
$ws.Range("D8").Value = 0.01226993865030675
$ws.Range("F8").Value = 0.04703476482617587
$ws.Range("J8").Value = 0.1226993865030675
$ws.Range("O8").Value = 0.016359918200409
$ws.Range("Q8").Value = 0.1738241308793456
$ws.Range("R8").Value = 0.1022494887525562
$ws.Range("S8").Value = 0.425357873210634
$ws.Range("B9").Value = 0.1165644171779141
$ws.Range("D9").Value = 0.03067484662576687
$ws.Range("F9").Value = 0.06134969325153374
$ws.Range("J9").Value = 0.1288343558282209
$ws.Range("O9").Value = 0.01226993865030675
$ws.Range("Q9").Value = 0.1901840490797546
$ws.Range("R9").Value = 0.09202453987730061
$ws.Range("S9").Value = 0.3680981595092024
$ws.Range("B10").Value = 0.1105730427764326
$ws.Range("D10").Value = 0.02340597255851493
$ws.Range("E10").Value = 0.001614205004035512
$ws.Range("F10").Value = 0.06376109765940274
$ws.Range("J10").Value = 0.122679580306699
$ws.Range("O10").Value = 0.01372074253430186
$ws.Range("Q10").Value = 0.1848264729620662
$ws.Range("R10").Value = 0.1089588377723971
$ws.Range("S10").Value = 0.3704600484261501
$ws.Range("G11").Value = 0.1372549019607843
$ws.Range("J11").Value = 0.1294117647058824
$ws.Range("K11").Value = 0.2313725490196079
$ws.Range("L11").Value = 0.4980392156862745
$ws.Range("S11").Value = 0.00392156862745098
$ws.Range("G12").Value = 0.7099236641221374
$ws.Range("J12").Value = 0.2442748091603053
$ws.Range("L12").Value = 0.02290076335877863
$ws.Range("S12").Value = 0.02290076335877863
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.01005025125628141
$ws.Range("H15").Value = 0.1809045226130653
$ws.Range("I15").Value = 0.08040201005025126
$ws.Range("J15").Value = 0.3316582914572864
$ws.Range("K15").Value = 0.06532663316582915
$ws.Range("M15").Value = 0.01005025125628141
$ws.Range("O15").Value = 0.04020100502512563
$ws.Range("S15").Value = 0.2814070351758794
$ws.Range("F16").Value = 0.005235602094240838
$ws.Range("H16").Value = 0.1727748691099476
$ws.Range("I16").Value = 0.1047120418848168
$ws.Range("J16").Value = 0.4659685863874345
$ws.Range("K16").Value = 0.07329842931937172
$ws.Range("M16").Value = 0.01570680628272251
$ws.Range("N16").Value = 0.005235602094240838
$ws.Range("O16").Value = 0.03141361256544502
$ws.Range("S16").Value = 0.1256544502617801
$ws.Range("F17").Value = 0.01
$ws.Range("H17").Value = 0.1975
$ws.Range("I17").Value = 0.0775
$ws.Range("J17").Value = 0.445
$ws.Range("K17").Value = 0.08500000000000001
$ws.Range("M17").Value = 0.0125
$ws.Range("O17").Value = 0.065
$ws.Range("S17").Value = 0.1075
$ws.Range("F18").Value = 0.01834862385321101
$ws.Range("H18").Value = 0.1880733944954129
$ws.Range("I18").Value = 0.06422018348623854
$ws.Range("J18").Value = 0.4954128440366973
$ws.Range("K18").Value = 0.0963302752293578
$ws.Range("M18").Value = 0.01376146788990826
$ws.Range("O18").Value = 0.04128440366972477
$ws.Range("S18").Value = 0.08256880733944955
$ws.Range("F19").Value = 0.009290540540540541
$ws.Range("H19").Value = 0.254222972972973
$ws.Range("I19").Value = 0.07010135135135136
$ws.Range("J19").Value = 0.3673986486486486
$ws.Range("K19").Value = 0.09375
$ws.Range("M19").Value = 0.02027027027027027
$ws.Range("N19").Value = 0.0008445945945945946
$ws.Range("O19").Value = 0.07601351351351351
$ws.Range("S19").Value = 0.1081081081081081
